$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the header formatting from F1 (style s="1")
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Update existing metric values in row 2
$ws.Range("B2").Value = 0.0564692933138092
$ws.Range("C2").Value = 0.99923544062592
$ws.Range("D2").Value = 0.1915777482460909

$pipelineText = "Pipeline(steps=[('model'," + [char]10 + "                 RandomForestRegressor(max_depth=5, n_estimators=150))])"
$ws.Range("F2").Value = $pipelineText

# Re-autofit the row so the embedded newline doesn't leave a stale custom row height
$ws.Rows(2).EntireRow.AutoFit()

# New elapsed time / cpu values
$ws.Range("G2").Value = 0.1289622459000384
$ws.Range("H2").Value = 0.991
